$d = $word.ActiveDocument

# Change 1: "Leaders and the management team..." paragraph text rewrite
$d.Content.Find.Execute(
    "hiring travel nurses and providing bonuses, and more. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "hiring travel nurses, and providing bonuses without losing current Nurses.",
    2
)
